# Formed the consolidated report
# Fill in the "Absent" column (H) so it correctly reflects attendance:
# Absent = 1 when the student was not present that day (Real/"E" column = 0),
# Absent = 0 when the student was present that day (Real/"E" column = 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 3) { $lastRow = 21 }

for ($row = 3; $row -le $lastRow; $row++) {
    $real = $ws.Cells.Item($row, 5).Value2  # column E = "Real"
    if ($real -eq 1) {
        $ws.Cells.Item($row, 8).Value = 0
    } else {
        $ws.Cells.Item($row, 8).Value = 1
    }
}
